# Updated CVDs for the month
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Manila Philippines - Internal Fill Rate / Commit-Forecast (row 5)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Manila Philippines")
$ws.Range("E5").Value = 0
$ws.Range("K5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 0
$ws.Range("U5").Value = 0
$ws.Range("V5").Value = 0
$ws.Range("W5").Value = 0

# ---------------------------------------------------------------------------
# Milwaukee Pmc Hq Wisconsin
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Milwaukee Pmc Hq Wisconsin")

# Professional Voluntary Turnover / Commit-Forecast (row 4)
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

# Internal Fill Rate / Commit-Forecast (row 7)
$ws.Range("O7").ClearContents()

# ---------------------------------------------------------------------------
# Monterrey Rbm Mexico
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Monterrey Rbm Mexico")

# Professional Voluntary Turnover: PY Actual / AOP / Commit-Forecast (rows 2-4)
$ws.Range("E2").Value = 0.1887
$ws.Range("E3").Value = 0.1887
$ws.Range("E4").Value = 0.1887

# Professional Voluntary Turnover / Commit-Forecast (row 4) monthly values
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 0
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0
$ws.Range("V4").Value = 0
$ws.Range("W4").Value = 0

# Manufacturing Voluntary Turnover: PY Actual / AOP / Commit-Forecast (rows 7-9)
$ws.Range("E7").Value = 0.3478
$ws.Range("E8").Value = 0.3478
$ws.Range("E9").Value = 0.3478

# Manufacturing Voluntary Turnover / Commit-Forecast (row 9) monthly values
$ws.Range("O9").Value = 0.0435
$ws.Range("P9").Value = 0.0496833333333333
$ws.Range("Q9").Value = 0.0496833333333333
$ws.Range("R9").Value = 0.14905
$ws.Range("S9").Value = 0.0496833333333333
$ws.Range("T9").Value = 0.0496833333333333
$ws.Range("U9").Value = 0.0496833333333333
$ws.Range("V9").Value = 0.14905
$ws.Range("W9").Value = 0.5962

# ---------------------------------------------------------------------------
# Rosemont Illinois - Internal Fill Rate / Commit-Forecast (row 7)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Rosemont Illinois")
$ws.Range("O7").ClearContents()

# ---------------------------------------------------------------------------
# Tipp City Ohio
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Tipp City Ohio")

# Professional Voluntary Turnover: PY Actual / AOP / Commit-Forecast (rows 2-4)
$ws.Range("E2").Value = 0.7143
$ws.Range("E3").Value = 0.7143
$ws.Range("E4").Value = 0.7143

# Professional Voluntary Turnover / Commit-Forecast (row 4) monthly values
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0.0833333333333333
$ws.Range("Q4").Value = 0.0833333333333333
$ws.Range("R4").Value = 0.25
$ws.Range("S4").Value = 0.0833333333333333
$ws.Range("T4").Value = 0.0833333333333333
$ws.Range("U4").Value = 0.0833333333333333
$ws.Range("V4").Value = 0.25
$ws.Range("W4").Value = 1

# ---------------------------------------------------------------------------
# Guadalajara Mexico - Professional Voluntary Turnover / Commit-Forecast (row 4)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Guadalajara Mexico")
$ws.Range("O4").ClearContents()

# ---------------------------------------------------------------------------
# Faridabad India - Professional Voluntary Turnover / Commit-Forecast (row 4)
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Faridabad India")
$ws.Range("O4").ClearContents()
